$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 1142.4286
$ws.Range("I106").Value = 969.4
$ws.Range("J106").Value = 1575
$ws.Range("K106").Value = 969.4
$ws.Range("L106").Value = 1575
$ws.Range("M106").Value = -338.4
$ws.Range("N106").Value = -2837

$ws.Range("H129").Value = 1025.9125
$ws.Range("J129").Value = 1039.909
$ws.Range("L129").Value = 3119.727
$ws.Range("N129").Value = -13119.727

$ws.Range("H135").Value = 2780149.2
$ws.Range("I135").Value = 3495
$ws.Range("K135").Value = 31455
$ws.Range("M135").Value = -28920

$ws.Range("H137").Value = 1046.6
$ws.Range("I137").Value = 958.0323
$ws.Range("J137").Value = 1242.7142
$ws.Range("K137").Value = 2874.0969
$ws.Range("L137").Value = 3728.1426
$ws.Range("M137").Value = -324.0969
$ws.Range("N137").Value = -8828.142599999999

$ws.Range("H138").Value = 4238.4673
$ws.Range("I138").Value = 2143.2856
$ws.Range("J138").Value = 5998.42
$ws.Range("K138").Value = 6429.8568
$ws.Range("L138").Value = 17995.26
$ws.Range("M138").Value = -1289.8568
$ws.Range("N138").Value = -28275.26

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8946.549000000001
$ws.Range("I32").Value = 7537.6206
$ws.Range("J32").Value = 29376
$ws.Range("K32").Value = 7537.6206
$ws.Range("L32").Value = 29376
$ws.Range("M32").Value = -7250.6206
$ws.Range("N32").Value = -29950

$ws.Range("H61").Value = 2173.2307
$ws.Range("I61").Value = 1960.2
$ws.Range("K61").Value = 1960.2
$ws.Range("M61").Value = -1748.2

$ws.Range("H102").Value = 2626.6667
$ws.Range("I102").Value = 1450
$ws.Range("J102").Value = 4980
$ws.Range("K102").Value = 1450
$ws.Range("L102").Value = 4980
$ws.Range("M102").Value = 172
$ws.Range("N102").Value = -8224

$ws.Range("H132").Value = 2123.8
$ws.Range("I132").Value = 1653.0238
$ws.Range("J132").Value = 2983.4783
$ws.Range("K132").Value = 4959.0714
$ws.Range("L132").Value = 8950.4349
$ws.Range("M132").Value = -2429.0714
$ws.Range("N132").Value = -14010.4349

$ws.Range("H136").Value = 2173.2307
$ws.Range("I136").Value = 1960.2
$ws.Range("K136").Value = 5880.6
$ws.Range("M136").Value = -3330.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H52").Value = 27873.428
$ws.Range("J52").Value = 27873.428
$ws.Range("L52").Value = 27873.428
$ws.Range("N52").Value = -28399.428

$ws.Range("H86").Value = 3100.625
$ws.Range("I86").Value = 3000.9092
$ws.Range("J86").Value = 3320
$ws.Range("K86").Value = 3000.9092
$ws.Range("L86").Value = 3320
$ws.Range("M86").Value = -1877.9092
$ws.Range("N86").Value = -5566

$ws.Range("H89").Value = 3100.625
$ws.Range("I89").Value = 3000.9092
$ws.Range("J89").Value = 3320
$ws.Range("K89").Value = 15004.546
$ws.Range("L89").Value = 16600
$ws.Range("M89").Value = -9388.546
$ws.Range("N89").Value = -27832

$ws.Range("H121").Value = 27873.428
$ws.Range("J121").Value = 27873.428
$ws.Range("L121").Value = 27873.428
$ws.Range("N121").Value = -31367.428

$ws.Range("H132").Value = 45726
$ws.Range("J132").Value = 45726
$ws.Range("L132").Value = 45726
$ws.Range("N132").Value = -55846

$ws.Range("H134").Value = 1858.878
$ws.Range("I134").Value = 1373.6786
$ws.Range("J134").Value = 2903.923
$ws.Range("K134").Value = 4121.0358
$ws.Range("L134").Value = 8711.769
$ws.Range("M134").Value = -1586.0358
$ws.Range("N134").Value = -13781.769

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3368174
$ws.Range("I22").Value = 4631101.5
$ws.Range("J22").Value = 366.66666
$ws.Range("K22").Value = 4631101.5
$ws.Range("L22").Value = 366.66666
$ws.Range("M22").Value = -4630751.5
$ws.Range("N22").Value = -1066.66666

$ws.Range("H31").Value = 4000.3052
$ws.Range("I31").Value = 2098.2239
$ws.Range("J31").Value = 8551.714
$ws.Range("K31").Value = 2098.2239
$ws.Range("L31").Value = 8551.714
$ws.Range("M31").Value = -1803.2239
$ws.Range("N31").Value = -9141.714

$ws.Range("H34").Value = 4000.3052
$ws.Range("I34").Value = 2098.2239
$ws.Range("J34").Value = 8551.714
$ws.Range("K34").Value = 2098.2239
$ws.Range("L34").Value = 8551.714
$ws.Range("M34").Value = -1896.2239
$ws.Range("N34").Value = -8955.714

$ws.Range("H132").Value = 835366.3
$ws.Range("I132").Value = 1250898.5
$ws.Range("J132").Value = 4301.9
$ws.Range("K132").Value = 3752695.5
$ws.Range("L132").Value = 12905.7
$ws.Range("M132").Value = -3750165.5
$ws.Range("N132").Value = -17965.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1901.6428
$ws.Range("I97").Value = 1980
$ws.Range("J97").Value = 1895.6154
$ws.Range("K97").Value = 5940
$ws.Range("L97").Value = 5686.8462
$ws.Range("M97").Value = -5444
$ws.Range("N97").Value = -6678.8462

$ws.Range("H131").Value = 773.73

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2416.25
$ws.Range("I126").Value = 2396.2856
$ws.Range("J126").Value = 2444.2
$ws.Range("K126").Value = 7188.8568
$ws.Range("L126").Value = 7332.599999999999
$ws.Range("M126").Value = -4718.8568
$ws.Range("N126").Value = -12272.6

$ws.Range("H132").Value = 2155.8958
$ws.Range("I132").Value = 1862.6
$ws.Range("J132").Value = 2945.5386
$ws.Range("K132").Value = 5587.799999999999
$ws.Range("L132").Value = 8836.6158
$ws.Range("M132").Value = -3057.799999999999
$ws.Range("N132").Value = -13896.6158

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3646.5
$ws.Range("I82").Value = 3544.75
$ws.Range("K82").Value = 3544.75
$ws.Range("M82").Value = -3183.75

$ws.Range("H85").Value = 3646.5
$ws.Range("I85").Value = 3544.75
$ws.Range("K85").Value = 3544.75
$ws.Range("M85").Value = -2296.75

$ws.Range("H132").Value = 10953.477
$ws.Range("I132").Value = 4045.0715
$ws.Range("J132").Value = 24770.285
$ws.Range("K132").Value = 12135.2145
$ws.Range("L132").Value = 74310.855
$ws.Range("M132").Value = -9605.2145
$ws.Range("N132").Value = -79370.855

$ws.Range("H136").Value = 3750.75
$ws.Range("I136").Value = 3825.1177
$ws.Range("K136").Value = 11475.3531
$ws.Range("M136").Value = -8925.3531

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 35697.5
$ws.Range("J109").Value = 35697.5
$ws.Range("L109").Value = 35697.5
$ws.Range("N109").Value = -38471.5

$ws.Range("H132").Value = 1338.7089
$ws.Range("I132").Value = 1065.6786
$ws.Range("J132").Value = 2003.4783
$ws.Range("K132").Value = 3197.0358
$ws.Range("L132").Value = 6010.4349
$ws.Range("M132").Value = -667.0357999999997
$ws.Range("N132").Value = -11070.4349
